# Generate Report for Handback
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" / "Error Detail" columns for row 7 (the b71050ac-... file) on both
# the zh-cn and de-de sheets, since a handback was processed for that file but
# it turned out to be based on a stale (non-latest) version of the source doc.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da05b0c79b9f223461a4b615da2323a12737c112/e2e/b71050ac-68ef-460f-a148-cea5c9c6805d.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9010361fbc699d2f03a2974b5a05be9947be90db/e2e/b71050ac-68ef-460f-a148-cea5c9c6805d.md."
$currentUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/da05b0c79b9f223461a4b615da2323a12737c112/e2e/b71050ac-68ef-460f-a148-cea5c9c6805d.md"
$targetDisplay = "b71050ac-68ef-460f-a148-cea5c9c6805d.md"

# Per-sheet data: language-specific handback xlf name + handback datetime
$sheetInfo = @{
    "zh-cn" = @{ Xlf = "b71050ac-68ef-460f-a148-cea5c9c6805d.3e0a1e8953fa08c7fb6d7950f3cd68138d6f1f01.zh-cn.xlf"; DateTime = "2016-08-26 04:41:27" }
    "de-de" = @{ Xlf = "b71050ac-68ef-460f-a148-cea5c9c6805d.3e0a1e8953fa08c7fb6d7950f3cd68138d6f1f01.de-de.xlf"; DateTime = "2016-08-26 04:41:33" }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $info = $sheetInfo[$sheetName]

    # I7 - Latest Target File: filled in with a hyperlink to the (stale) current version
    $ws.Range("I7").Value = $targetDisplay
    $ws.Hyperlinks.Add($ws.Range("I7"), $currentUrl, "", "", $targetDisplay) | Out-Null

    # J7 - Latest Handback File
    $ws.Range("J7").Value = $info.Xlf

    # K7 - Latest Handback DateTime
    $ws.Range("K7").Value = $info.DateTime

    # P7 - Error Detail
    $ws.Range("P7").Value = $errorDetail

    # Column P (16) is now wide enough to show the long error message.
    $ws.Columns.Item(16).ColumnWidth = 39.17
}
